$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CCDeferredPlanCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:50:40 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:52:10 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 22:20:02 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 22:21:16 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 22:01:17 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 22:02:35 IST 2025"

$ws = $wb.Worksheets.Item("CCDeferredPlanCredit")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:53:38 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:56:13 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:43:11 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:44:01 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:44:57 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:45:52 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:46:43 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:47:44 IST 2025"

$ws = $wb.Worksheets.Item("CCAutoPayPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:48:41 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:49:41 IST 2025"

$ws = $wb.Worksheets.Item("VerifyReceiptPageDataCC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 23:17:44 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 23:18:26 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:13:33 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:14:14 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 20:59:08 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:00:31 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:22:21 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:24:45 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:16:18 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:16:58 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:31:55 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:34:25 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:04:20 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:05:39 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 24 19:17:49 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 24 19:19:02 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredCredit")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:17:34 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:20:01 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredCredit")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 20:56:33 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 20:57:53 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelDeferredPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 21:14:55 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 21:15:35 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelDeferredPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:27:07 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:29:32 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelDeferredPC")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 21:01:48 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 21:03:03 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanPS")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 22:22:34 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 22:25:01 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCredit")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 22:15:16 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 22:17:40 IST 2025"

$ws = $wb.Worksheets.Item("CMCDeferredPlanCorp")
$ws.Range("A2").Value = "Fail"
$ws.Range("B2").Value = "Tue Jun 03 22:10:32 IST 2025"
$ws.Range("A3").Value = "Fail"
$ws.Range("B3").Value = "Tue Jun 03 22:12:54 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 22:03:55 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 22:04:45 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanCredit")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 22:05:32 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 22:06:27 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 22:07:14 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 22:08:00 IST 2025"

$ws = $wb.Worksheets.Item("CMCAutoPayPlanPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 22:08:55 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 22:09:46 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:06:46 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:07:32 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:05:32 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:06:09 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:08:08 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:08:49 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:12:35 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:13:22 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:10:48 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:11:47 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:15:59 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:16:51 IST 2025"

$ws = $wb.Worksheets.Item("VerifyEditLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:14:13 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:15:08 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCorp")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:00:13 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:01:01 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayCC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 22:58:37 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 22:59:25 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPC")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:01:54 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:02:49 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCancelLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:03:44 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:04:37 IST 2025"

$ws = $wb.Worksheets.Item("VerifyCreateLabelAutopayPS")
$ws.Range("A2").Value = "Pass"
$ws.Range("B2").Value = "Tue Jun 03 23:09:30 IST 2025"
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Tue Jun 03 23:10:07 IST 2025"
